$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header in H1, matching the style (bold, border, centered) used
# by the other header cells (e.g. G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New "Save" column values for rows 2..7
$ws.Cells.Item(2, 8).Value = 1
$ws.Cells.Item(3, 8).Value = 0
$ws.Cells.Item(4, 8).Value = 1
$ws.Cells.Item(5, 8).Value = 0
$ws.Cells.Item(6, 8).Value = 1
$ws.Cells.Item(7, 8).Value = 0
